$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to be treated as text so numeric-looking strings
# (e.g. "313.20", "1.002") are not auto-converted to numbers,
# matching the inline-string cell type used in the source file.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '28.362.99'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '1.808.19'
$ws.Range('E3').Value = '  -0.84%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = '313.20'
$ws.Range('E5').Value = '  -0.92%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('D7').Value = '0.5157'
$ws.Range('E7').Value = '  -0.38%  '
$ws.Range('D8').Value = '0.4004'
$ws.Range('E8').Value = '  +3.74%  '
$ws.Range('D9').Value = '0.07872'
$ws.Range('E9').Value = '  -4.88%  '
$ws.Range('D10').Value = '1.112'
$ws.Range('E10').Value = '  -0.84%  '
$ws.Range('D11').Value = '40.90'
$ws.Range('E11').Value = '  -2.33%  '
$ws.Range('D12').Value = '6.356'
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('D13').Value = '1.002'
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('D14').Value = '20.42'
$ws.Range('E14').Value = '  -3.26%  '
$ws.Range('D15').Value = '7.316'
$ws.Range('E15').Value = '  -2.15%  '
$ws.Range('D16').Value = '1.806.60'
$ws.Range('E16').Value = '  -1.14%  '
$ws.Range('D17').Value = '92.65'
$ws.Range('E17').Value = '  -1.34%  '
$ws.Range('D18').Value = '0.00001087'
$ws.Range('E18').Value = '  -2.92%  '
$ws.Range('D19').Value = '0.06575'
$ws.Range('E19').Value = '  -0.87%  '
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  -0.24%  '
$ws.Range('D21').Value = '17.31'
$ws.Range('E21').Value = '  -2.80%  '
$ws.Range('D22').Value = '6.004'
$ws.Range('E22').Value = '  -0.79%  '
$ws.Range('D23').Value = '28.385.10'
$ws.Range('E23').Value = '  -0.42%  '
$ws.Range('D24').Value = '11.13'
$ws.Range('E24').Value = '  -3.03%  '
$ws.Range('D25').Value = '2.227'
$ws.Range('E25').Value = '  -0.90%  '
$ws.Range('D26').Value = '160.82'
$ws.Range('E26').Value = '  +1.00%  '
$ws.Range('D27').Value = '20.54'
$ws.Range('E27').Value = '  -2.41%  '
$ws.Range('D28').Value = '2.021.82'
$ws.Range('E28').Value = '  -0.72%  '
$ws.Range('D29').Value = '2.408'
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').Value = '128.26'
$ws.Range('E30').Value = '  +1.85%  '
$ws.Range('D31').Value = '0.1102'
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('D32').Value = '1.066'
$ws.Range('E32').Value = '  -2.57%  '
$ws.Range('D33').Value = '3.670'
$ws.Range('E33').Value = '  -0.44%  '
$ws.Range('D34').Value = '5.572'
$ws.Range('E34').Value = '  -2.68%  '
$ws.Range('D35').Value = '0.07210'
$ws.Range('E35').Value = '  -4.88%  '
$ws.Range('D36').Value = '9.122'
$ws.Range('E36').Value = '  +4.32%  '
$ws.Range('D37').Value = '0.02343'
$ws.Range('E37').Value = '  -1.05%  '
$ws.Range('D38').Value = '0.2181'
$ws.Range('E38').Value = '  -1.96%  '
$ws.Range('D39').Value = '11.59'
$ws.Range('E39').Value = '  -4.19%  '
$ws.Range('D40').Value = '5.052'
$ws.Range('E40').Value = '  -3.67%  '
$ws.Range('D41').Value = '0.6199'
$ws.Range('E41').Value = '  -3.14%  '
$ws.Range('D42').Value = '1.001'
$ws.Range('E42').Value = '  -0.22%  '
$ws.Range('D43').Value = '1.155'
$ws.Range('E43').Value = '  -3.01%  '
$ws.Range('D44').Value = '13.25'
$ws.Range('E44').Value = '  -2.97%  '
$ws.Range('D45').Value = '0.5996'
$ws.Range('E45').Value = '  -3.61%  '
$ws.Range('D46').Value = '1.304'
$ws.Range('E46').Value = '  -6.71%  '
$ws.Range('D47').Value = '3.740'
$ws.Range('E47').Value = '  -1.58%  '
$ws.Range('D48').Value = '125.52'
$ws.Range('E48').Value = '  -1.84%  '
$ws.Range('D49').Value = '1.219'
$ws.Range('E49').Value = '  +1.20%  '
$ws.Range('D50').Value = '1.926'
$ws.Range('E50').Value = '  -4.00%  '
$ws.Range('D51').Value = '0.06831'
$ws.Range('E51').Value = '  -1.94%  '

# Reset the style back to Normal so no stray s="n" attribute is
# left on the cells (keeps them matching the original, unstyled cells).
$ws.Range("D2:E51").Style = "Normal"
